$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new computed metric values
$ws.Range("B2").Value = 4.3429963328524117
$ws.Range("C2").Value = 0.36932979120863829
$ws.Range("D2").Value = 3.462684439839427
$ws.Range("E2").Value = 0.34322378219884236
$ws.Range("F2").Value = 0.58585303805548572
$ws.Range("G2").Value = 0.53031679417311295
$ws.Range("H2").Value = 0.65677621780115758
$ws.Range("I2").Value = 0.83454519166440044

$ws.Range("B3").Value = 4.6116412802035098
$ws.Range("C3").Value = 0.39217544308356667
$ws.Range("D3").Value = 3.5799294346158832
$ws.Range("E3").Value = 0.38699867984222958
$ws.Range("F3").Value = 0.62209217953791185
$ws.Range("G3").Value = 0.5482730910413246
$ws.Range("H3").Value = 0.61300132015777042
$ws.Range("I3").Value = 0.79707356012971065

$ws.Range("B4").Value = 4.0798707592104977
$ws.Range("C4").Value = 0.34695350863162883
$ws.Range("D4").Value = 3.218397337635309
$ws.Range("E4").Value = 0.30289441057403549
$ws.Range("F4").Value = 0.55035843826912978
$ws.Range("G4").Value = 0.49290375375620021
$ws.Range("H4").Value = 0.69710558942596457
$ws.Range("I4").Value = 0.84361841737650511

# Add new row 5 for the old_model metrics
$ws.Range("A5").Value = "old_model"
$ws.Range("B5").Value = 4.0628397039906616
$ws.Range("C5").Value = 0.34550518227205451
$ws.Range("D5").Value = 3.0211847857374399
$ws.Range("E5").Value = 0.30037087752345837
$ws.Range("F5").Value = 0.54806101624131087
$ws.Range("G5").Value = 0.46270027142616538
$ws.Range("H5").Value = 0.69962912247654163
$ws.Range("I5").Value = 0.84531908412890899

# Widen column A slightly to fit the new "old_model" label
$ws.Columns.Item(1).ColumnWidth = 14.45
